$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New poll rows appended to the bottom of the data table:
#  - row 149: opinionway poll
#  - row 150: opinionway poll
#  - row 151: cluster17 poll

$ws.Cells.Item(149, 1).Value  = 55
$ws.Cells.Item(149, 2).Value  = 2021
$ws.Cells.Item(149, 3).Value  = 17
$ws.Cells.Item(149, 4).Value  = 12
$ws.Cells.Item(149, 5).Value  = 22
$ws.Cells.Item(149, 6).Value  = "opinionway"
$ws.Cells.Item(149, 7).Value  = "online"
$ws.Cells.Item(149, 8).Value  = "included"
$ws.Cells.Item(149, 9).Value  = 985
$ws.Cells.Item(149, 10).Value = 1
$ws.Cells.Item(149, 11).Value = "T_1"
$ws.Cells.Item(149, 12).Value = 8
$ws.Cells.Item(149, 13).Value = 2
$ws.Cells.Item(149, 14).Value = 1
$ws.Cells.Item(149, 15).Value = 6
$ws.Cells.Item(149, 16).Value = 4
$ws.Cells.Item(149, 17).Value = 26
$ws.Cells.Item(149, 18).Value = 18
$ws.Cells.Item(149, 21).Value = 1
$ws.Cells.Item(149, 22).Value = 2
$ws.Cells.Item(149, 23).Value = 16
$ws.Cells.Item(149, 24).Value = 12
$ws.Cells.Item(149, 29).Value = 3

$ws.Cells.Item(150, 1).Value  = 56
$ws.Cells.Item(150, 2).Value  = 2022
$ws.Cells.Item(150, 3).Value  = 18
$ws.Cells.Item(150, 4).Value  = 1
$ws.Cells.Item(150, 5).Value  = 2
$ws.Cells.Item(150, 6).Value  = "opinionway"
$ws.Cells.Item(150, 7).Value  = "online"
$ws.Cells.Item(150, 8).Value  = "included"
$ws.Cells.Item(150, 9).Value  = 1059
$ws.Cells.Item(150, 10).Value = 1
$ws.Cells.Item(150, 11).Value = 1
$ws.Cells.Item(150, 12).Value = 9
$ws.Cells.Item(150, 13).Value = 3
$ws.Cells.Item(150, 14).Value = 1
$ws.Cells.Item(150, 15).Value = 7
$ws.Cells.Item(150, 16).Value = 4
$ws.Cells.Item(150, 17).Value = 26
$ws.Cells.Item(150, 18).Value = 16
$ws.Cells.Item(150, 21).Value = 1
$ws.Cells.Item(150, 22).Value = 2
$ws.Cells.Item(150, 23).Value = 16
$ws.Cells.Item(150, 24).Value = 13

$ws.Cells.Item(151, 1).Value  = 57
$ws.Cells.Item(151, 2).Value  = 2021
$ws.Cells.Item(151, 3).Value  = 18
$ws.Cells.Item(151, 4).Value  = 12
$ws.Cells.Item(151, 5).Value  = 29
$ws.Cells.Item(151, 6).Value  = "cluster17"
$ws.Cells.Item(151, 7).Value  = "online"
$ws.Cells.Item(151, 8).Value  = "partially"
$ws.Cells.Item(151, 9).Value  = 2176
$ws.Cells.Item(151, 10).Value = 1
$ws.Cells.Item(151, 11).Value = 0.5
$ws.Cells.Item(151, 12).Value = 13
$ws.Cells.Item(151, 13).Value = 1.5
$ws.Cells.Item(151, 14).Value = 1.5
$ws.Cells.Item(151, 15).Value = 4
$ws.Cells.Item(151, 16).Value = 2
$ws.Cells.Item(151, 17).Value = 23
$ws.Cells.Item(151, 18).Value = 15
$ws.Cells.Item(151, 21).Value = 1
$ws.Cells.Item(151, 22).Value = 1.5
$ws.Cells.Item(151, 23).Value = 14.5
$ws.Cells.Item(151, 24).Value = 15
$ws.Cells.Item(151, 25).Value = 1
$ws.Cells.Item(151, 27).Value = 1
$ws.Cells.Item(151, 29).Value = 4.5

# Move the selection to the new last cell, matching the view state recorded
# after the edit.
$ws.Range("AD151").Select()
